$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ---------------------------------------------------------------------------
# 1) Reverse the order of the "Periodo Mora" column (E16:E104) so periods run
#    ascending (oldest -> newest) instead of descending. 89 data rows.
# ---------------------------------------------------------------------------
$firstDataRow = 16
$lastDataRow = 104

$periods = @()
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
  $periods += $ws.Range("E" + $r).Value()
}

$count = $periods.Length
for ($i = 0; $i -lt $count; $i++) {
  $r = $firstDataRow + $i
  $val = $periods[$count - 1 - $i]
  $ws.Range("E" + $r).Value = $val
}

# ---------------------------------------------------------------------------
# 2) Insert a brand-new data row right after the old last row (105) for the
#    newest period "2508" (part 1 of the new statement). Everything below
#    (the signature block) shifts down by one row automatically.
# ---------------------------------------------------------------------------
$newRow = $lastDataRow + 1
$ws.Rows.Item($newRow).Insert()

$ws.Range("B" + $newRow).Value = "CC"
$ws.Range("C" + $newRow).Value = "72146160"
$ws.Range("D" + $newRow).Value = "HENRY ALFONSO ATENCIO MONTAÑO"
$ws.Range("E" + $newRow).Value = "2508"
$ws.Range("F" + $newRow).Value = 200000
$ws.Range("G" + $newRow).Value = 5000000

# Match the borders/format of the neighbouring data rows so the new row and
# the (no-longer-last) previous row both look consistent.
$dataRange = $ws.Range("B" + $firstDataRow + ":J" + $newRow)

$ws.Range("B" + $newRow + ":J" + $newRow).Borders.Item(9).LineStyle = 1
$ws.Range("B" + $newRow + ":J" + $newRow).Borders.Item(9).Weight = 2
$ws.Range("B" + $newRow + ":J" + $newRow).Borders.Item(7).LineStyle = 1
$ws.Range("B" + $newRow + ":J" + $newRow).Borders.Item(10).LineStyle = 1
$ws.Range("B" + $newRow + ":J" + $newRow).Borders.Item(8).LineStyle = 1

$ws.Range("F" + $newRow + ":G" + $newRow).NumberFormat = $ws.Range("F" + $lastDataRow).NumberFormat()

# ---------------------------------------------------------------------------
# 3) Update the summary figures at the top of the statement.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 18000000
$ws.Range("F13").Value = 90
